## EPEXSPOT prices workbook - daily automated data refresh
## Adds the 17-jul column to "Prix Spot" and the 2025-07-15 row to
## "Gaz" and "CO2".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Prix Spot": append column AH ("17-jul") after AG ("16-jul")
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Prix Spot")

$ws.Range("AH1").Value = "17-jul"

$ahValues = @{
    2  = 104.23
    3  = 98.38
    4  = 90.1
    5  = 83.77
    6  = 82.12
    7  = 89.16
    8  = 87.05
    9  = 101.62
    10 = 105.92
    11 = 93.1
    12 = 85
    13 = 77.07
    14 = 70.17
    15 = 56.14
    16 = 51.21
    17 = 66.56
    18 = 76.81
    19 = 84.36
    20 = 91.52
    21 = 108.66
    22 = 121.41
    23 = 125.8
    24 = 125.07
    25 = 114.18
}

foreach ($row in $ahValues.Keys) {
    $ws.Cells.Item($row, 34).Value = $ahValues[$row]
}

# Match the look of the rest of the header row (bold, centered, bordered)
[void]$ws.Range("AG1").Copy()
[void]$ws.Range("AH1").PasteSpecial(-4122)
[void]$ws.Range("A1").Select()

# ---------------------------------------------------------------
# Sheet "Gaz": append row 31 (2025-07-15, 33.35)
# Dates in column A are stored as plain text, not real dates, so we
# force text formatting before the write (and drop it again right
# after) to stop the auto date-recognition from turning the literal
# "2025-07-15" into a date serial number.
# ---------------------------------------------------------------
$gaz = $wb.Worksheets.Item("Gaz")
$gaz.Range("A31").NumberFormat = "@"
$gaz.Range("A31").Value = "2025-07-15"
$gaz.Range("A31").ClearFormats()
$gaz.Range("B31").Value = 33.35

# ---------------------------------------------------------------
# Sheet "CO2": append row 31 (2025-07-15, 70.8)
# ---------------------------------------------------------------
$co2 = $wb.Worksheets.Item("CO2")
$co2.Range("A31").NumberFormat = "@"
$co2.Range("A31").Value = "2025-07-15"
$co2.Range("A31").ClearFormats()
$co2.Range("B31").Value = 70.8

Write-Host "Prix Spot / Gaz / CO2 updated with 2025-07-15 data"
